# Update CDA Logical model for ST.r2b
# - Bump Version / Date metadata values
# - Insert a new "Jurisdiction" property row (empty value) right after "Contact"
#   on the Metadata sheet (sheet1); everything below shifts down one row.
# - The Elements sheet (sheet2) references the same shared strings and will
#   automatically follow along (its own cell content is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Insert a new row for "Jurisdiction" right after the "Contact" row (row 10),
#    copying formatting from the row above so the style matches the rest of the table.
$ws.Rows.Item(11).Insert()
$ws.Range("A10:B10").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# 2. Update the Version value (now on row 3, unaffected by the insert above).
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# 3. Update the Date value (now on row 8, unaffected by the insert above).
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"
